$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values are plain decimals (e.g. "231.50")
# that Excel would otherwise auto-convert to a number, silently dropping
# the trailing zero / exact text (e.g. "0.310" -> 0.31). Mark just those
# specific cells as Text up front so the literal string is preserved,
# matching every other price cell on the sheet, which is stored as text.
$textRefs = @("D5", "D8", "D9", "D14", "D16", "D18", "D20", "D21", "D22", "D25", "D27", "D29", "D32", "D36", "D37", "D38", "D39", "D41", "D47", "D48", "D51")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "35.471.75"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").Value = "1.838.75"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "231.50"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "43.77"
$ws.Range("E8").Value = "  +12.71%  "
$ws.Range("D9").Value = "0.310"
$ws.Range("E9").Value = "  +7.48%  "
$ws.Range("E10").Value = "  +5.24%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "2.104.49"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.835.28"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "0.674"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "4.71"
$ws.Range("E16").Value = "  +7.75%  "
$ws.Range("D17").Value = "35.401.69"
$ws.Range("D18").Value = "70.12"
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  +4.37%  "
$ws.Range("D20").Value = "244.42"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "12.12"
$ws.Range("E21").Value = "  +8.48%  "
$ws.Range("D22").Value = "4.70"
$ws.Range("E22").Value = "  +14.95%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").Value = "170.85"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "17.78"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").Value = "1.61"
$ws.Range("E29").Value = "  +31.38%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Value = "3.332.69"
$ws.Range("E31").Value = "  +37.17%  "
$ws.Range("D32").Value = "0.0554"
$ws.Range("E32").Value = "  +7.59%  "
$ws.Range("E33").Value = "  +6.19%  "
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").Value = "96.02"
$ws.Range("E36").Value = "  +16.62%  "
$ws.Range("D37").Value = "0.689"
$ws.Range("E37").Value = "  +7.33%  "
$ws.Range("D38").Value = "1.12"
$ws.Range("E38").Value = "  +6.16%  "
$ws.Range("D39").Value = "15.59"
$ws.Range("E39").Value = "  +11.40%  "
$ws.Range("D40").Value = "1.349.79"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").Value = "2.45"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("E43").Value = "  +6.11%  "
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "6.27"
$ws.Range("E47").Value = "  +8.48%  "
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "2.007.34"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "103.41"
$ws.Range("E51").Value = "  +0.70%  "
